$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.164.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.566.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.566.18"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.174.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.566.92"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.296.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.578"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.709.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +26.29%  "
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.568.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.09"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0817"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  +5.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.439.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  +0.22%  "
